# (Integration Test) Room Detail
# Fills in the "Room Detail" test-case rows (8-11) that were previously
# blank placeholders, adds a matching "back to previous page" test case
# row (28) on "House Detail", and restores the saved sheet-view cursor
# positions recorded for both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# House Detail (sheet 3): row 28 - "back to previous of page" test case
# ---------------------------------------------------------------------
$house = $wb.Worksheets.Item("House Detail")

$house.Range("A28").Value = "TS_FHF_HouseDetail_021"
$house.Range("B28").Value = "Verify the back to previous of page functionality of House Detail"
$house.Range("C28").Value = "TC_FHF_HouseDetail_Back_001"
$house.Range("D28").Value = "Click link on top of page"
$house.Range("E28").Value = "1. Click link on top of page"
$house.Range("F28").Value = "Valid URL`nTest Data"
$house.Range("H28").Value = "User back to previous page"
$house.Range("I28").Value = "System backs to previous page"
$house.Range("J28").Value = "System backs to previous page"
$house.Range("K28").Value = "Pass"

# ---------------------------------------------------------------------
# Room Detail (sheet 4): rows 8-11 - new integration test cases
# ---------------------------------------------------------------------
$room = $wb.Worksheets.Item("Room Detail")

# Row 8 - show room images
$room.Range("A8").Value = "TS_FHF_RoomDetail_001"
$room.Range("B8").Value = "Verify the show image of room functionality of Room Detail"
$room.Range("C8").Value = "TC_FHF_RoomDetail_ShowRoomImage_001"
$room.Range("D8").Value = "Load page"
$room.Range("E8").Value = "1. Load page"
$room.Range("F8").Value = "Valid URL`nTest Data"
$room.Range("H8").Value = "User see images of room"
$room.Range("I8").Value = "System shows the images of room"
$room.Range("J8").Value = "System shows the images of room"
$room.Range("K8").Value = "Pass"
$room.Range("L8").Value = "KienNT"

# Row 9 - show room information
$room.Range("A9").Value = "TS_FHF_RoomDetail_002"
$room.Range("B9").Value = "Verify the show information of room functionality of Room Detail"
$room.Range("C9").Value = "TC_FHF_RoomDetail_ShowRoomInformation_001"
$room.Range("D9").Value = "Load page"
$room.Range("E9").Value = "1. Load page"
$room.Range("F9").Value = "Valid URL`nTest Data"
$room.Range("H9").Value = "User see information of room"
$room.Range("I9").Value = "System shows the information of room"
$room.Range("J9").Value = "System shows the information of room"
$room.Range("K9").Value = "Pass"
$room.Range("L9").Value = "KienNT"

# Row 10 - show landlord detail
$room.Range("A10").Value = "TS_FHF_RoomDetail_003"
$room.Range("B10").Value = "Verify the detail landlord of house functionality of Room Detail"
$room.Range("C10").Value = "TC_FHF_RoomDetail_ShowLandlord_001"
$room.Range("D10").Value = "Load page"
$room.Range("E10").Value = "1. Load page"
$room.Range("F10").Value = "Valid URL`nTest Data"
$room.Range("H10").Value = "User see detail information of landlord "
$room.Range("I10").Value = "System shows the detail information of landlord"
$room.Range("J10").Value = "System shows the detail information of landlord"
$room.Range("K10").Value = "Pass"
$room.Range("L10").Value = "KienNT"

# Row 11 - back to previous page
$room.Range("A11").Value = "TS_FHF_HouseDetail_004"
$room.Range("B11").Value = "Verify the back to previous of page functionality of Room Detail"
$room.Range("C11").Value = "TC_FHF_RoomDetail_Back_001"
$room.Range("D11").Value = "Click link on top of page"
$room.Range("E11").Value = "1. Click link on top of page"
$room.Range("F11").Value = "Valid URL`nTest Data"
$room.Range("H11").Value = "User back to previous page"
$room.Range("I11").Value = "System backs to previous page"
$room.Range("J11").Value = "System backs to previous page"
$room.Range("K11").Value = "Pass"
$room.Range("L11").Value = "KienNT"

# ---------------------------------------------------------------------
# Restore the saved cursor / scroll position for each sheet view
# ---------------------------------------------------------------------
$houseWindow = $excel.ActiveWindow
$house.Activate()
$excel.ActiveWindow.ScrollRow = 19
$house.Range("C27").Select()

$room.Activate()
$excel.ActiveWindow.Zoom = 70
$room.Range("J14").Select()

$house.Activate()
